# Fruta / hortaliza, semanal
# Insert 7 new rows of Mandarina price data at row 612 (pushing existing
# rows 612-673 down to 619-680), and populate the newly inserted rows
# with the new weekly data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 blank rows starting at row 612; this shifts the existing
# rows 612:673 down to 619:680 and keeps their contents untouched.
$ws.Range("A612:A618").EntireRow.Insert()

# Constant columns shared by every data row in this block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102004
$categoria = "Mandarina"

# Data for the 7 newly inserted rows (612-618).
$rows = @(
    @{ Row=612; Fecha=44826; Variedad="Murcott"; Calidad="Especial"; Volumen=480;  PMin=8000; PMax=8000; PProm=8000; Unidad="$/bandeja 10 kilos"; Origen="Provincia de Limarí";          PKg=800; KgUnidad=10 },
    @{ Row=613; Fecha=44826; Variedad="Murcott"; Calidad="Especial"; Volumen=240;  PMin=7000; PMax=7000; PProm=7000; Unidad="$/bandeja 10 kilos"; Origen="Región Metropolitana";         PKg=700; KgUnidad=10 },
    @{ Row=614; Fecha=44826; Variedad="Murcott"; Calidad="Primera";  Volumen=600;  PMin=7000; PMax=7000; PProm=7000; Unidad="$/bandeja 10 kilos"; Origen="Provincia de Limarí";          PKg=700; KgUnidad=10 },
    @{ Row=615; Fecha=44826; Variedad="Murcott"; Calidad="Primera";  Volumen=390;  PMin=6000; PMax=6000; PProm=6000; Unidad="$/bandeja 10 kilos"; Origen="Región Metropolitana";         PKg=600; KgUnidad=10 },
    @{ Row=616; Fecha=44826; Variedad="Murcott"; Calidad="Segunda";  Volumen=360;  PMin=6000; PMax=6000; PProm=6000; Unidad="$/bandeja 10 kilos"; Origen="Provincia de Limarí";          PKg=600; KgUnidad=10 },
    @{ Row=617; Fecha=44826; Variedad="Murcott"; Calidad="Segunda";  Volumen=240;  PMin=5000; PMax=5000; PProm=5000; Unidad="$/bandeja 10 kilos"; Origen="Región Metropolitana";         PKg=500; KgUnidad=10 },
    @{ Row=618; Fecha=44826; Variedad="Murcott"; Calidad="Tercera";  Volumen=240;  PMin=4000; PMax=4000; PProm=4000; Unidad="$/bandeja 10 kilos"; Origen="Región Metropolitana";         PKg=400; KgUnidad=10 }
)

foreach ($r in $rows) {
    $ri = $r.Row
    $ws.Cells.Item($ri, 1).Value  = $mercadoId
    $ws.Cells.Item($ri, 2).Value  = $mercado
    $ws.Cells.Item($ri, 3).Value  = $region
    $ws.Cells.Item($ri, 4).Value  = $r.Fecha
    $ws.Cells.Item($ri, 5).Value  = $codreg
    $ws.Cells.Item($ri, 6).Value  = $tipo
    $ws.Cells.Item($ri, 7).Value  = $productoId
    $ws.Cells.Item($ri, 8).Value  = $producto
    $ws.Cells.Item($ri, 9).Value  = $categoriaId
    $ws.Cells.Item($ri, 10).Value = $categoria
    $ws.Cells.Item($ri, 11).Value = $r.Variedad
    $ws.Cells.Item($ri, 12).Value = $r.Calidad
    $ws.Cells.Item($ri, 13).Value = $r.Volumen
    $ws.Cells.Item($ri, 14).Value = $r.PMin
    $ws.Cells.Item($ri, 15).Value = $r.PMax
    $ws.Cells.Item($ri, 16).Value = $r.PProm
    $ws.Cells.Item($ri, 17).Value = $r.Unidad
    $ws.Cells.Item($ri, 18).Value = $r.Origen
    $ws.Cells.Item($ri, 19).Value = $r.PKg
    $ws.Cells.Item($ri, 20).Value = $r.KgUnidad
}
